# Update "想去人数" (want-to-go count) figures in column F across the
# four sheets of the workbook, reflecting refreshed scrape numbers.
# (gh-pages data refresh @ 456a3b4)

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" --------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 12947
$ws.Range("F3").Value  = 7262
$ws.Range("F4").Value  = 266
$ws.Range("F8").Value  = 162
$ws.Range("F9").Value  = 374
$ws.Range("F10").Value = 1059
$ws.Range("F13").Value = 1037
$ws.Range("F15").Value = 283
$ws.Range("F16").Value = 387
$ws.Range("F18").Value = 294
$ws.Range("F19").Value = 326
$ws.Range("F21").Value = 272
$ws.Range("F22").Value = 409
$ws.Range("F23").Value = 5333
$ws.Range("F24").Value = 81
$ws.Range("F25").Value = 1481
$ws.Range("F27").Value = 2979
$ws.Range("F29").Value = 75
$ws.Range("F30").Value = 1405
$ws.Range("F35").Value = 3760

# --- Sheet "演出" --------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value  = 3747
$ws.Range("F3").Value  = 3747
$ws.Range("F6").Value  = 86
$ws.Range("F18").Value = 50

# --- Sheet "本地生活" -----------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9330
$ws.Range("F3").Value = 571
$ws.Range("F4").Value = 2079

# --- Sheet "全部类型" -----------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 9330
$ws.Range("F3").Value  = 571
$ws.Range("F4").Value  = 2079
$ws.Range("F5").Value  = 12947
$ws.Range("F6").Value  = 7262
$ws.Range("F7").Value  = 3747
$ws.Range("F9").Value  = 162
$ws.Range("F10").Value = 374
$ws.Range("F11").Value = 1059
$ws.Range("F14").Value = 1037
$ws.Range("F16").Value = 283
$ws.Range("F17").Value = 387
$ws.Range("F19").Value = 294
$ws.Range("F20").Value = 326
$ws.Range("F25").Value = 272
$ws.Range("F26").Value = 409
$ws.Range("F27").Value = 5333
$ws.Range("F28").Value = 81
$ws.Range("F29").Value = 1481
$ws.Range("F34").Value = 2984
$ws.Range("F36").Value = 75
$ws.Range("F37").Value = 1405
$ws.Range("F47").Value = 3760
$ws.Range("F49").Value = 50

$wb.Save()
